$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update several odds values ---
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.75
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 5
$ws.Range("X2").Value = 8
$ws.Range("AE2").Value = 21
$ws.Range("AK2").Value = 51
$ws.Range("AS2").Value = 301
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 26
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151

# --- Row 4: update several odds values ---
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 2
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.05
$ws.Range("V4").Value = 1.7
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 19
$ws.Range("AG4").Value = 501
$ws.Range("AN4").Value = 5.5
$ws.Range("AS4").Value = 351
$ws.Range("AT4").Value = 2.38
$ws.Range("AY4").Value = 26

# --- Rows 5/6/7: the Millonarios vs Dep. Pasto fixture (old row 5) was
# removed from the sheet, so the remaining fixtures shift up one row:
# old row 6 -> new row 5, old row 7 -> new row 6. Capture the data first
# (arrays), then delete the now-superseded last row.
$oldRow6 = $ws.Range("A6:BD6").Value2
$oldRow7 = $ws.Range("A7:BD7").Value2

$ws.Range("A5:BD5").Value2 = $oldRow6
$ws.Range("A6:BD6").Value2 = $oldRow7

# Remove the now-duplicated trailing row and let Excel shrink the sheet
# dimension (A1:BD7 -> A1:BD6) automatically.
$ws.Range("A7:BD7").EntireRow.Delete()
